$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit swaps the full data content of row 12 and row 13 (two sighting
# records that were reordered), including which optional columns are
# populated vs. left blank for each record.
#
# Row 12 currently holds the "Garnlav" (lichen) record; row 13 currently
# holds the "Järpe" (bird) record. After the edit, row 12 holds the bird
# record and row 13 holds the lichen record.

# Helper: paste a literal value into a cell as *text*, even if it looks
# numeric, without leaving any NumberFormat/style change behind on the
# destination cell. Uses an off-sheet scratch cell + TEXT() + copy/paste
# values (mirrors how a user would paste text-from-formula results).
function Set-TextValue($rangeRef, [string]$text) {
    $scratch = $ws.Range("BZ1000")
    $scratch.Formula = '=TEXT("' + $text + '","@")'
    $scratch.Copy()
    $ws.Range($rangeRef).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# --- Simple numeric value swaps: both rows already have a value
#     present in these columns, so a straight swap is enough ---
$swapCols = @("A", "B", "E", "Q", "R")
foreach ($col in $swapCols) {
    $c12 = $ws.Range($col + "12")
    $c13 = $ws.Range($col + "13")
    $v12 = $c12.Value2
    $v13 = $c13.Value2
    $c12.Value2 = $v13
    $c13.Value2 = $v12
}

# --- Text value swaps (species name / scientific name / author) ---
$textSwapCols = @("F", "G", "H")
foreach ($col in $textSwapCols) {
    $c12 = $ws.Range($col + "12")
    $c13 = $ws.Range($col + "13")
    $v12 = [string]$c12.Value2
    $v13 = [string]$c13.Value2
    $c12.Value = $v13
    $c13.Value = $v12
}

# --- Antal (I): row 12 blank -> "2", row 13 "2" -> blank (both stay
#     present; "2" must remain text, not be coerced to a number) ---
Set-TextValue "I12" "2"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = ""

# --- Metod (N): row 12 blank -> "observerad", row 13 "observerad" -> blank
#     (both stay present) ---
$ws.Range("N12").Value = "observerad"
$ws.Range("N13").NumberFormat = "@"
$ws.Range("N13").Value = ""

# --- Enhet (J): present-but-blank on row 12, entirely absent on row 13
#     -> row 12 becomes absent, row 13 becomes present-but-blank ---
$ws.Range("J12").ClearContents()
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = ""

# --- Bestämningsmetod (AF): present-but-blank on row 12, entirely absent
#     on row 13 -> row 12 becomes absent, row 13 becomes present-but-blank ---
$ws.Range("AF12").ClearContents()
$ws.Range("AF13").NumberFormat = "@"
$ws.Range("AF13").Value = ""

# --- Kön (L): entirely absent on row 12, present-but-blank on row 13
#     -> row 12 becomes present-but-blank, row 13 becomes absent ---
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = ""
$ws.Range("L13").ClearContents()

# --- Aktivitet (M): absent on row 12, "födosökande" on row 13
#     -> row 12 gets the text, row 13 becomes absent ---
$ws.Range("M12").Value = "födosökande"
$ws.Range("M13").ClearContents()

# --- Publik kommentar (AC): absent on row 12, comment text on row 13
#     -> row 12 gets the text, row 13 becomes absent ---
$ws.Range("AC12").Value = "Synobservation av 2 st födosökande järpar."
$ws.Range("AC13").ClearContents()
